$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the additional "new title" test data column.
# (Single-quoted so PowerShell does not try to interpolate "${new_title}"
#  as a variable reference.)
$ws.Range("C1").Value = '${new_title}'

# New rows of test-case data appended below the existing ones
$ws.Range("A5").Value = "empty image"
$ws.Range("A6").Value = "removing image"
$ws.Range("B6").Value = "remove_image"
$ws.Range("A7").Value = "test case 1"
$ws.Range("B7").Value = "same_title"
$ws.Range("A8").Value = "sample editing test case"
$ws.Range("B8").Value = "edit_new_title"
$ws.Range("C8").Value = "new title"

# Widen the columns to fit the new, longer content
$ws.Columns.Item(1).ColumnWidth = 22.833333333333336
$ws.Columns.Item(2).ColumnWidth = 20.333333333333332
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666

# Update the active selection to match the author's final cursor position
$ws.Range("B8").Select()
